$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ticket count cells to be descriptive text values ("N Boletos")
# instead of bare numbers. This causes them to become shared-string cells.
$ws.Range("B2").Value = "6 Boletos"
$ws.Range("B3").Value = "7 Boletos"

# Maximize the workbook window (Excel was reopened/maximized before saving).
$excel.ActiveWindow.WindowState = -4137

# Move/update the current selection to E7, matching the last selection
# recorded when the file was saved.
$ws.Range("E7").Select()
